$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97 (pushes old rows 97..145 down to 98..146)
$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with the new data record
$ws.Range("A97").Value2 = 7
$ws.Range("B97").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C97").Value2 = "Ñuble"
$ws.Range("D97").Value2 = 44460
$ws.Range("E97").Value2 = 16
$ws.Range("F97").Value2 = 100112043
$ws.Range("G97").Value2 = "Pepino ensalada"
$ws.Range("H97").Value2 = "Sin especificar"
$ws.Range("I97").Value2 = "Primera"
$ws.Range("J97").Value2 = 160
$ws.Range("K97").Value2 = 16000
$ws.Range("L97").Value2 = 17000
$ws.Range("M97").Value2 = 16500
$ws.Range("N97").Value2 = "$/caja 60 unidades"
$ws.Range("O97").Value2 = "Región del Maule"
$ws.Range("P97").Value2 = 275
$ws.Range("Q97").Value2 = 60
$ws.Range("R97").Value2 = "Hortaliza"
